{"js": "// Remove the last row of the \"Actores y Funciones\" table \u2014 the\n// \"Proveedor\" / \"Venta de mercanc\u00eda a Reposter\u00eda\" row.\nconst table = context.document.body.tables.getFirst();\nconst rows = table.rows;\nrows.load(\"items\");\nawait context.sync();\n\nconst lastRow = rows.items[rows.items.length - 1];\nlastRow.delete();\nawait context.sync();\n", "ps1": "# Remove the last row of the \"Actores y Funciones\" table \u2014 the\n# \"Proveedor\" / \"Venta de mercanc\u00eda a Reposter\u00eda\" row.\n$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n$lastRowIndex = $t.Rows.Count\n$t.Rows.Item($lastRowIndex).Delete()\n"}
